# Updates crypto price/volume/name/link values per the authoritative diff.
# For numeric-looking "Price" strings (e.g. "225.85"), a leading apostrophe is used
# (the standard Excel text-entry prefix) so the cell stays a text value, matching
# the original inlineStr cell type, instead of being auto-converted to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.492.94"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.809.17"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'225.85"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("D9").Value = "'0.292"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "'0.0968"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").Value = "2.071.01"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "1.818.50"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "34.469.63"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "'68.41"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "'242.31"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").Value = "0.0₃0775"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("E24").Value = "  +4.95%  "
$ws.Range("D25").Value = "'171.48"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "'7.88"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "'17.23"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").Value = "'0.122"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.23"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'3.91"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("D34").Value = "'1.79"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "1.363.58"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("E36").Value = "  -4.27%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = "  -6.79%  "
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("D42").Value = "'80.86"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "'0.936"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").Value = "'1.17"
$ws.Range("E44").Value = "  +4.78%  "
$ws.Range("D45").Value = "'13.38"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").Value = "1.971.66"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "'5.82"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'102.40"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  -6.54%  "
